# Fixed a bug in stats: rows in the data table were being written out of
# order relative to their id/reel values. Re-apply the corrected row order
# (same rows, same per-row values, just placed on the correct output row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected data, in the row order it should appear in (rows 2..25).
# Columns: A (symbol id), B (reel1), C (reel2), D (reel3), E (reel4), F (reel5)
$data = @(
    @(101,  9, 30, 15, 60, 15),
    @(201,  9, 30, 15, 45, 30),
    @(1203, 3, 15, 15, 15, 15),
    @(501,  9, 52, 30, 75, 45),
    @(801,  3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(901, 16, 15, 45, 60, 60),
    @(902,  1,  0,  0,  0,  0),
    @(301,  6, 45, 30, 60, 45),
    @(401,  9, 48, 67, 75, 45),
    @(1001,18, 30, 75, 60, 72),
    @(701,  3, 90, 45, 97, 15),
    @(601,  9, 60, 67, 60, 42),
    @(1,    0,  2,  2,  2,  2),
    @(2,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(1101, 0, 15, 30, 30,  0),
    @(3,    0,  3,  3,  3,  3),
    @(802,  0,  4,  5,  4,  0),
    @(602,  0,  0,  4,  0,  9),
    @(402,  0,  0,  4,  0,  0),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
